# MCA_labels.xlsx - "Cleaning up and more model training"
#
# This script:
#  1. Flips a number of 0 -> 1 label cells in rows 49-70 (columns O..W range,
#     differs per row) to reflect newly-labelled training examples.
#  2. Highlights column A (the row-index column) for rows 49-70 with the
#     same light-green fill used for the existing rows (style index 1 /
#     fill FFE2EFDA), matching the "cleaning up" pass over old rows.
#  3. Appends a brand-new row 71 (index value 70) for one more labelled
#     training example, with columns U-X marked 1 and everything else 0,
#     and the same green highlight on its A cell.
#  4. Updates the sheet view / selection to mirror where the user ended up
#     after adding the new row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Light green fill already used by style index 1 in this workbook
# (fill FFE2EFDA -> OLE BGR value 14348258)
$highlightColor = 14348258

# ---------------------------------------------------------------------------
# 1) Flip existing 0 -> 1 cells for rows 49-70 (contiguous ranges per row)
# ---------------------------------------------------------------------------
$flips = @{
    49 = "T49:V49"
    50 = "S50:V50"
    51 = "S51:U51"
    52 = "R52:U52"
    53 = "S53:V53"
    54 = "R54:U54"
    55 = "T55:V55"
    56 = "R56:T56"
    57 = "T57:U57"
    58 = "P58:R58"
    59 = "S59:V59"
    60 = "R60:T60"
    61 = "T61:U61"
    62 = "Q62:T62"
    63 = "O63:S63"
    64 = "S64:V64"
    65 = "R65:T65"
    66 = "S66:U66"
    67 = "Q67:T67"
    68 = "Q68:U68"
    69 = "R69:T69"
    70 = "T70:W70"
}

for ($r = 49; $r -le 70; $r++) {
    $ws.Range($flips[$r]).Value = 1
    # Highlight the row-index cell the same way the rest of the sheet does
    $ws.Range("A$r").Interior.Color = $highlightColor
}

# ---------------------------------------------------------------------------
# 2) Append new row 71 for a newly-labelled example
# ---------------------------------------------------------------------------
$ws.Range("A71:AW71").Value = 0
$ws.Range("A71").Value = 70
$ws.Range("U71:X71").Value = 1
$ws.Range("A71").Interior.Color = $highlightColor

# ---------------------------------------------------------------------------
# 3) Update view state to where the editor ended up (best effort - some
#    view-state attributes such as topLeftCell are not persisted by this
#    runtime, but the selection is).
# ---------------------------------------------------------------------------
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 61
$ws.Range("A70:A71").Select()

Write-Host "Applied MCA_labels edits: flipped labels rows 49-70, added row 71, highlighted A49:A71"
